$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newQuery = "SELECT `n    DISTINCT cf.file_name AS ""File Name"",`n    CASE`n        WHEN cf.file_name LIKE '%.bai' THEN 'bai'`n        WHEN cf.file_name LIKE '%.bam' THEN 'bam'`n        WHEN cf.file_name LIKE '%.csv' THEN 'csv'`n        WHEN cf.file_name LIKE '%.doc' THEN 'doc'`n        WHEN cf.file_name LIKE '%.docx' THEN 'docx'`n           WHEN cf.file_name LIKE '%.fastq.gz' THEN 'fastq'`n        WHEN cf.file_name LIKE '%.vcf.gz' THEN 'vcf'`n        WHEN cf.file_name LIKE '%.tsv.gz' THEN 'tsv'`n        WHEN cf.file_name LIKE '%.pdf' THEN 'pdf'`n        WHEN cf.file_name LIKE '%.rtf' THEN 'rtf'`n        WHEN cf.file_name LIKE '%.tbi' THEN 'tbi'`n        WHEN cf.file_name LIKE '%.tif' THEN 'tif'`n        WHEN cf.file_name LIKE '%.xls' THEN 'xls'`n        WHEN cf.file_name LIKE '%.xlsx' THEN 'xlsx'`n        ELSE 'Unknown'`n    END AS ""Format"",`n    cf.file_type AS ""File Type"",`n    CASE     `n    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN `n        CASE `n            WHEN ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'`n            ELSE ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'`n        END`n    WHEN cf.file_size >= 1024 * 1024 THEN `n        CASE `n            WHEN ROUND(cf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'`n            ELSE ROUND(cf.file_size / (1024.0 * 1024.0), 2) || ' MB'`n        END`n    WHEN cf.file_size >= 1024 THEN `n        CASE `n            WHEN ROUND(cf.file_size / 1024.0, 2) = CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'`n            ELSE ROUND(cf.file_size / 1024.0, 2) || ' KB'`n        END`n    ELSE `n        CASE `n            WHEN ROUND(cf.file_size, 2) = CAST(ROUND(cf.file_size, 0) AS INT) `n            THEN CAST(CAST(ROUND(cf.file_size, 0) AS INT) AS TEXT) || ' Bytes'`n            ELSE ROUND(cf.file_size, 2) || ' Bytes'`n        END`nEND AS ""Size"",`n    'sample' AS ""Association"",`n    cf.file_description AS ""Description"",`n    smp.sample_id AS ""Sample ID"",`n    c.case_record_id AS ""Case ID"",`n    dmg.breed AS ""Breed"",`n    diag.disease_term AS ""Diagnosis""`nFROM `n    df_case_file cf`nJOIN `n    df_sample smp ON cf.""sample.sample_id"" = smp.sample_id`nJOIN `n    df_case c ON smp.""case.case_record_id"" = c.case_record_id`nJOIN `n    df_study st ON c.""study.clinical_study_designation"" = st.clinical_study_designation`nJOIN `n    df_program p ON st.""program.program_acronym"" = p.program_acronym`nJOIN `n    df_demographic dmg ON dmg.""case.case_record_id"" = c.case_record_id`nJOIN `n    df_diagnosis diag ON diag.""case.case_record_id"" = c.case_record_id`nJOIN `n    df_enrollment enr ON enr.""case.case_record_id"" = c.case_record_id`nJOIN `n    df_publication pub ON pub.""study.clinical_study_designation"" = st.clinical_study_designation`nLEFT JOIN `n    df_study_file sf ON sf.""study.clinical_study_designation"" = st.clinical_study_designation`nWHERE`n    st.clinical_study_designation = 'COTC022' AND cf.file_type = 'RNA Sequence File'`nORDER BY `n    cf.file_name ASC`nLIMIT 100;"

$ws.Range("B4").Value2 = $newQuery

# Update the selection / view state on the sheet to match the saved file (row 5 active)
$ws.Range("A5").Select()
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("B5").Select()
